# Insert a new data row at row 216, shifting all existing rows (216-280)
# down by one (to 217-281), then populate the new row 216 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 216. Excel will push the
# existing row 216 (and everything below it) down to row 217, extending
# the sheet from A1:R280 to A1:R281, and the new row inherits formatting
# (e.g. the date style on column D) from the row above it.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with its data.
$row = 216
$ws.Cells.Item($row, 1).Value  = 6
$ws.Cells.Item($row, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($row, 3).Value  = "Metropolitana"
$ws.Cells.Item($row, 4).Value  = 44782
$ws.Cells.Item($row, 5).Value  = 13
$ws.Cells.Item($row, 6).Value  = 100112026
$ws.Cells.Item($row, 7).Value  = "Haba"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 500
$ws.Cells.Item($row, 11).Value = 12000
$ws.Cells.Item($row, 12).Value = 14000
$ws.Cells.Item($row, 13).Value = 13080
$ws.Cells.Item($row, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Coquimbo"
$ws.Cells.Item($row, 16).Value = 523
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
